# Move "SEMANA DE PROVAS" exam week two rows earlier (to 22-Mar and 27-Mar),
# and shift the intervening content (CSP lecture, then the adversarial-search
# lecture repeated once more) down accordingly. The date column (A) is left
# untouched; only columns B-E (Questão/Fundamentos/Conteúdo/Programação) move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (22-Mar) currently holds the CSP lecture content; capture it before
# it gets overwritten, it needs to move down to row 15 (29-Mar).
$cspB = $ws.Cells.Item(13, 2).Value2
$cspC = $ws.Cells.Item(13, 3).Value2
$cspD = $ws.Cells.Item(13, 4).Value2
$cspE = $ws.Cells.Item(13, 5).Value2

# Row 17 (05-Abr) holds the "adversarial search" lecture content, which also
# needs to be duplicated onto row 16 (03-Abr), the row vacated by the exam.
$advB = $ws.Cells.Item(17, 2).Value2
$advC = $ws.Cells.Item(17, 3).Value2
$advD = $ws.Cells.Item(17, 4).Value2
$advE = $ws.Cells.Item(17, 5).Value2

# Row 15 (29-Mar) currently holds the "SEMANA DE PROVAS" exam-week content.
$provasB = $ws.Cells.Item(15, 2).Value2
$provasC = $ws.Cells.Item(15, 3).Value2
$provasD = $ws.Cells.Item(15, 4).Value2
$provasE = $ws.Cells.Item(15, 5).Value2

# Rows 13 and 14 (22-Mar, 27-Mar) become the exam week.
$ws.Cells.Item(13, 2).Value = $provasB
$ws.Cells.Item(13, 3).Value = $provasC
$ws.Cells.Item(13, 4).Value = $provasD
$ws.Cells.Item(13, 5).Value = $provasE

$ws.Cells.Item(14, 2).Value = $provasB
$ws.Cells.Item(14, 3).Value = $provasC
$ws.Cells.Item(14, 4).Value = $provasD
$ws.Cells.Item(14, 5).Value = $provasE

# Row 15 (29-Mar) now gets the CSP lecture that used to be taught on 22-Mar.
$ws.Cells.Item(15, 2).Value = $cspB
$ws.Cells.Item(15, 3).Value = $cspC
$ws.Cells.Item(15, 4).Value = $cspD
$ws.Cells.Item(15, 5).Value = $cspE

# Row 16 (03-Abr) now gets the adversarial-search lecture content (same as
# rows 17/18).
$ws.Cells.Item(16, 2).Value = $advB
$ws.Cells.Item(16, 3).Value = $advC
$ws.Cells.Item(16, 4).Value = $advD
$ws.Cells.Item(16, 5).Value = $advE

$wb.Save()
